# Update workbook metadata (IG regeneration: new URL, new date, new
# "Jurisdiction" property row inserted after "Contact").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# URL changed from the "pythia" IG to the "cicada" IG.
$ws.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/valid-age-status"

# Regeneration timestamp.
$ws.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# Insert a new "Jurisdiction" property row right after "Contact" (row 10),
# pushing Description/Purpose/.../Context down by one row.
$ws.Rows.Item(11).Insert()

# Match the formatting used by the other property rows.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
